# Apply the "add fr jurisdiction, snomed parameters" edit to the workbook.
#
# Semantic changes (the diff is purely a shared-string table reshuffle except
# for these two content edits, both on the "Metadata" sheet):
#   - B11 (the "Jurisdiction" value cell) changes from "" to "FRANCE"
#   - B8  (the "Date" value cell) changes its timestamp

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B11").Value = "FRANCE"
$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"
